# Commit: "add 2d act camera"
# Update the default (villageScene) row's camera offset position/rotation
# so the scene uses a more "2D act" style camera.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 = villageScene (ID 1) -> J:CamOffestPos, K:CamOffestRot
$ws.Range("J2").Value = "0,4.2,5.5"
$ws.Range("K2").Value = "25,180"
